# TrialsSetup update 2026-02-27 12:00
# The "ALPINE" trial's "Days remaining" figure (row 8, column B) needs to be
# corrected from 13 to 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B8").Value = 12
